# mapping.xlsx - added table schemas for t_magic_item and related enums
#
# 1) Rename "magic_item_table" -> "t_magic_item" and extend its columns
#    with d_from/d_to (split from d_last_update) plus a PK/FK annotation row.
# 2) Add two new sheets "e_rarity" and "e_category" describing enum tables
#    (copied from t_magic_item's formatting/page setup), each with a
#    PK/FK annotation row.
# 3) Restore view state (selections / active tab / zoom) to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 4: magic_item_table -> t_magic_item
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "t_magic_item"

# Row 1 (columns): split "d_last_update" into "d_from" (I) / "d_to" (J),
# and move "t_write" from J to K.
$ws4.Cells.Item(1,9).Value  = "d_from"
$ws4.Cells.Item(1,10).Value = "d_to"
$ws4.Cells.Item(1,11).Value = "t_write"

# Row 2 (types): I/J/K are all timestamps.
$ws4.Cells.Item(2,9).Value  = "timestamp"
$ws4.Cells.Item(2,10).Value = "timestamp"
$ws4.Cells.Item(2,11).Value = "timestamp"

# Row 3 (new): PK / FK annotations.
$ws4.Cells.Item(3,1).Value = "PK"
$ws4.Cells.Item(3,5).Value = "FK e_category"
$ws4.Cells.Item(3,6).Value = "FK e_rarity"

# ---------------------------------------------------------------------
# Sheet 5 (new): e_rarity -- duplicate t_magic_item's sheet formatting
# (page setup, header/footer, etc.) then replace the content.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4.Copy([System.Reflection.Missing]::Value, $lastSheet)
$ws5 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5.Name = "e_rarity"

$ws5.Cells.Clear()
$ws5.Cells.EntireColumn.ClearFormats()

$ws5.Cells.Item(1,1).Value = "n_id"
$ws5.Cells.Item(1,2).Value = "s_name"
$ws5.Cells.Item(1,3).Value = "n_from_character_level"
$ws5.Cells.Item(1,4).Value = "n_value_from"
$ws5.Cells.Item(1,5).Value = "n_value_to"
$ws5.Cells.Item(1,6).Value = "b_valid"
$ws5.Cells.Item(1,7).Value = "t_write"

$ws5.Cells.Item(2,1).Value = "integer"
$ws5.Cells.Item(2,2).Value = "string"
$ws5.Cells.Item(2,3).Value = "integer"
$ws5.Cells.Item(2,4).Value = "integer"
$ws5.Cells.Item(2,5).Value = "integer"
$ws5.Cells.Item(2,6).Value = "boolean"
$ws5.Cells.Item(2,7).Value = "timestamp"

$ws5.Cells.Item(3,1).Value = "PK"

# column widths auto-fit around the longer header text
$ws5.Columns.Item(3).ColumnWidth = 20.0066666666667
$ws5.Columns.Item(4).ColumnWidth = 11.8066666666667

# ---------------------------------------------------------------------
# Sheet 6 (new): e_category -- same approach as e_rarity.
# ---------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4.Copy([System.Reflection.Missing]::Value, $lastSheet2)
$ws6 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6.Name = "e_category"

$ws6.Cells.Clear()
$ws6.Cells.EntireColumn.ClearFormats()

$ws6.Cells.Item(1,1).Value = "n_id"
$ws6.Cells.Item(1,2).Value = "s_name"
$ws6.Cells.Item(1,3).Value = "b_valid"
$ws6.Cells.Item(1,4).Value = "t_write"

$ws6.Cells.Item(2,1).Value = "integer"
$ws6.Cells.Item(2,2).Value = "string"
$ws6.Cells.Item(2,3).Value = "boolean"
$ws6.Cells.Item(2,4).Value = "timestamp"

$ws6.Cells.Item(3,1).Value = "PK"

# ---------------------------------------------------------------------
# View state: selections on each sheet, zoom on the new sheets, and the
# workbook's final active tab/sheet.
# ---------------------------------------------------------------------
$win = $wb.Windows.Item(1)

$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("B2").Select()

$ws4.Activate()
$ws4.Range("A4").Select()

$ws5.Activate()
$ws5.Range("A4").Select()
$win.Zoom = 110

$ws6.Activate()
$ws6.Range("D3").Select()
$win.Zoom = 110
